$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("three")

$data = @(
    @("Marhta", "Ankunding", "walrus", "(901) 409-9620", "Western West Virginia Academy", "Albany"),
    @("Quinn", "Schmitt", "hippopotamus", "(716) 580-6045", "Western Corwin College", "Nelson Bay"),
    @("Booker", "Friesen", "cricket", "(810) 310-3311 x8287", "North Olson University", "Cairns"),
    @("Becki", "Harber", "ferret", "(510) 320-9211 x0295", "East Bartell Academy", "Newcastle"),
    @("Gabriel", "Grimes", "ape", "(480) 971-5729 x7175", "Gerhold Academy", "Mount Gambier")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowValues[$j]
    }
}
